# Apply the mortality_2020 table updates (non-gaussian ranovas refresh).
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Table 1 caption (PVE values) ---
Replace-Text "PVE for population: 11.684. PVE for family: 8.19" "PVE for population: 15.818. PVE for family: 8.435"

# --- Table 1 body values ---
Replace-Text "0.1695" "0.0955"

# 0.037 -> 0.0575, and this run drops its bold formatting
$r1 = $d.Content
$r1.Find.Execute("0.037") | Out-Null
$r1.Text = "0.0575"
$r1.Font.Bold = $false

# --- Table 2 caption (PVE values) ---
Replace-Text "PVE for population: 11.398. PVE for family: 7.666" "PVE for population: 15.507. PVE for family: 7.919"

# --- Table 2 body values ---
Replace-Text "0.1635" "0.0905"
Replace-Text "0.052" "0.077"

# --- Table 3 (chi-squared table) header row height: 637 -> 571 twips (31.85 -> 28.55 pt) ---
$t3 = $d.Tables.Item(3)
$t3.Rows.Item(1).Height = 28.55

# --- Table 3: fix mangled chi symbol and update numeric values ---
Replace-Text "Ï‡" "χ"
Replace-Text "5.973" "6.132"
Replace-Text "0.113" "0.105"
Replace-Text "1.088" "1.173"
Replace-Text "0.297" "0.279"

# --- Table 4 caption (PVE values) ---
Replace-Text "PVE for population: 11.507. PVE for family: 7.892" "PVE for population: 15.624. PVE for family: 8.141"

# --- Table 4 grid: widen the third column 961 -> 1084 twips (48.05 -> 54.2 pt) ---
$t4 = $d.Tables.Item(4)
$t4.Columns.Item(3).Width = 54.2

# --- Table 4 body values ---
Replace-Text "0.166" "0.0925"

# 0.044 -> 0.0665, and this run drops its bold formatting
$r2 = $d.Content
$r2.Find.Execute("0.044") | Out-Null
$r2.Text = "0.0665"
$r2.Font.Bold = $false

# --- Table 5 (chi-squared table) header row height: 637 -> 571 twips (31.85 -> 28.55 pt) ---
$t5 = $d.Tables.Item(5)
$t5.Rows.Item(1).Height = 28.55

# --- Table 5: numeric value updates (chi symbol already fixed above, both occurrences) ---
Replace-Text "6.059" "6.226"
Replace-Text "0.109" "0.101"
Replace-Text "0.811" "0.864"
Replace-Text "0.368" "0.353"
